$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the NODE CLASS category labels used in column C (rows 2-17)
# "Gate Station" -> "Gas pumping station"
# "Regulator Station" -> "Intermediate station"
# "Other" -> "Deliver station"
$range = $ws.Range("C2:C17")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    $current = $cell.Value2
    if ($current -eq "Gate Station") { $cell.Value2 = "Gas pumping station" }
    elseif ($current -eq "Regulator Station") { $cell.Value2 = "Intermediate station" }
    elseif ($current -eq "Other") { $cell.Value2 = "Deliver station" }
}

# Update the current selection to match the saved workbook state
$ws.Range("D13").Select()
